# The sentence currently reads "... submit it with you program file ..."
# and needs the missing "r" restored so it reads "... with your program
# file ...". The author apparently selected a point after "you", typed
# "r", and then retyped the remainder of the sentence, which is why the
# final document has the single original run split into three runs that
# all share identical run formatting.

$d = $word.ActiveDocument

$old = ", save it into a .doc or .docx file and submit it with you program file as stated above"

# Locate the run's range precisely using Find (keeps us from hand
# computing character offsets that could drift).
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence."
}

$runStart = $rng.Start
$runEnd = $rng.End

$part1 = ", save it into a .doc or .docx file and submit it with you"
$part2 = "r"
$part3 = " program file as stated above"

$insertPos = $runStart + $part1.Length

# Insert the missing "r" so the text reads "... with your program ...".
$d.Range($insertPos, $insertPos).Text = $part2

# After the insertion, everything from $insertPos onward shifted right by
# the length of $part2 (1 character).
$p1Start = $runStart
$p1End = $insertPos
$p2Start = $insertPos
$p2End = $insertPos + $part2.Length
$p3Start = $p2End
$p3End = $runEnd + $part2.Length

# Sanity-check the three pieces now read back exactly as expected.
if ($d.Range($p1Start, $p1End).Text -ne $part1) { throw "part1 mismatch" }
if ($d.Range($p2Start, $p2End).Text -ne $part2) { throw "part2 mismatch" }
if ($d.Range($p3Start, $p3End).Text -ne $part3) { throw "part3 mismatch" }

# The edit above leaves the touched paragraph as one run with uniform
# formatting. Re-introduce the three run boundaries (matching how Word
# itself keeps freshly-typed text in its own run) by toggling a
# formatting property on/off across each piece, which forces the writer
# to split runs at those exact boundaries while leaving the effective
# formatting unchanged. Go right-to-left so earlier offsets stay valid.
$d.Range($p3Start, $p3End).Bold = 1
$d.Range($p3Start, $p3End).Bold = 0
$d.Range($p2Start, $p2End).Bold = 1
$d.Range($p2Start, $p2End).Bold = 0
$d.Range($p1Start, $p1End).Bold = 1
$d.Range($p1Start, $p1End).Bold = 0
